$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text relabeling (same cells, renamed headers). "year" is
#     written last so it lands at the end of the shared-string table, just
#     like the authored edit. ---
$ws.Range("B1").Value = "applications"
$ws.Range("C1").Value = "approved"
$ws.Range("D1").Value = "dollars mn"
$ws.Range("E1").Value = "dollars nyc"
$ws.Range("F1").Value = "dollars other"
$ws.Range("G1").Value = "genre-Visual Arts"
$ws.Range("H1").Value = "genre-Theater"
$ws.Range("I1").Value = "genre-Dance"
$ws.Range("J1").Value = "genre-Media Arts/Film and Video"
$ws.Range("K1").Value = "genre-Multidisciplinary"
$ws.Range("L1").Value = "genre-Literature"
$ws.Range("M1").Value = "genre-Music"
$ws.Range("N1").Value = "genre-Arts Criticism"
$ws.Range("O1").Value = "genre-Other Disciplines"
$ws.Range("A1").Value = "year"

# --- Header row is taller now ---
$ws.Rows(1).RowHeight = 80

# --- New number format (plain thousands separator, no $) for the three
#     dollar-amount data columns across all data rows ---
$ws.Range("D2:F51").NumberFormat = "#,##0"

# --- Selection moved to the newly (re)formatted data block ---
$ws.Range("D2:F51").Select()
